# Append: 2025-12-10 01:24 JST
# - A new job listing row is inserted at row 11 (pushing the former rows
#   11-20 down to 12-21).
# - Every listing's "取得日時" (fetched-at) timestamp is refreshed to the
#   new run time.
# - A fresh hyperlink relationship is registered for the row that rolled
#   off the bottom of the previously-seen range (old row 20 -> new row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-10 01:24:21"

# Insert a new row above the current row 11; Excel shifts rows 11-20 down
# to 12-21 and carries their formatting (incl. the F-column hyperlink
# style) along with them.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the freshly scraped listing.
$ws.Range("A11").Value = $newTimestamp
$ws.Range("B11").Value = "【フルスタックエンジニア募集】新規Webサービス開発"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5450548"
$ws.Range("G11").Value = 75
$ws.Range("H11").Value = "◆開発"

# Refresh the "取得日時" timestamp on every other already-present listing
# row (rows 2-10 stayed in place; rows that used to be 11-20 are now
# 12-21 after the insert above) - i.e. every data row except the brand
# new row 11, which was already stamped above.
for ($r = 2; $r -le 21; $r++) {
    if ($r -ne 11) {
        $ws.Cells.Item($r, 1).Value = $newTimestamp
    }
}

# The listing that rolled past the bottom of the previously tracked range
# (now row 21) gets a brand-new hyperlink relationship registered for it.
$ws.Hyperlinks.Add($ws.Range("F21"), "https://www.lancers.jp/work/detail/5449948")

# Match the existing URL-column look (the other F-cells already carry the
# workbook's "Hyperlink" cell style from the row-11 insert above).
$ws.Range("F21").Style = $ws.Range("F20").Style
